$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.7408638000488281
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 5332.369331643707
$ws.Range("F2").Value = 0.1224743717815145
$ws.Range("G2").Value = 0.1119448212796044
$ws.Range("H2").Value = 0.1119448212796044
$ws.Range("I2").Value = 0.1119448212796044
$ws.Range("J2").Value = 0.1119448212796044
$ws.Range("K2").Value = 0.1119448212796044
$ws.Range("L2").Value = 0.1119448212796044
$ws.Range("M2").Value = 0.1119448212796044
$ws.Range("N2").Value = 0.1119448212796044
$ws.Range("O2").Value = 0.1119448212796044
$ws.Range("P2").Value = 0.1119448212796044
$ws.Range("Q2").Value = 0.1119448212796044
$ws.Range("R2").Value = 0.1119448212796044
$ws.Range("S2").Value = 0.1119448212796044
$ws.Range("T2").Value = 0.1119448212796044
$ws.Range("U2").Value = 0.1119448212796044
$ws.Range("V2").Value = 0.1119448212796044
$ws.Range("W2").Value = 0.1119448212796044
$ws.Range("X2").Value = 0.1119448212796044
$ws.Range("Y2").Value = 0.1119448212796044

$ws.Range("C3").Value = 0.8230006694793701
$ws.Range("E3").Value = 2283.877863671712
$ws.Range("F3").Value = 0.1320051911585163
$ws.Range("G3").Value = 0.125711923168711
$ws.Range("H3").Value = 0.1007998690413547
$ws.Range("I3").Value = 0.08365886483087463
$ws.Range("J3").Value = 0.0749109377536249
$ws.Range("K3").Value = 0.06883631511461699
$ws.Range("L3").Value = 0.06270459288546704
$ws.Range("M3").Value = 0.05935795837570595
$ws.Range("N3").Value = 0.05606197634175334
$ws.Range("O3").Value = 0.05108948225248054
$ws.Range("P3").Value = 0.05062690877992004
$ws.Range("Q3").Value = 0.049519271491518
$ws.Range("R3").Value = 0.04790161215254951
$ws.Range("S3").Value = 0.04703083164544249
$ws.Range("T3").Value = 0.04592965306132922
$ws.Range("U3").Value = 0.04547553417312739
$ws.Range("V3").Value = 0.04524014485222019
$ws.Range("W3").Value = 0.04493049377152694
$ws.Range("X3").Value = 0.04460052186854468
$ws.Range("Y3").Value = 0.04452003632888327

$ws.Range("C4").Value = 0.734485387802124
$ws.Range("E4").Value = 2360.759664688449
$ws.Range("G4").Value = 0.1228077057349246
$ws.Range("H4").Value = 0.1028632017443405
$ws.Range("I4").Value = 0.09046948163690653
$ws.Range("J4").Value = 0.07903592409496828
$ws.Range("K4").Value = 0.06760585506952446
$ws.Range("L4").Value = 0.06399615381852641
$ws.Range("M4").Value = 0.06103205789511927
$ws.Range("N4").Value = 0.05734910240102561
$ws.Range("O4").Value = 0.05554845181227069
$ws.Range("P4").Value = 0.05317126230851055
$ws.Range("Q4").Value = 0.05165127808603515
$ws.Range("R4").Value = 0.05059626683613808
$ws.Range("S4").Value = 0.0494474970517505
$ws.Range("T4").Value = 0.04850820813741588
$ws.Range("U4").Value = 0.04736591095766077
$ws.Range("V4").Value = 0.04720031978818517
$ws.Range("W4").Value = 0.04662258797467365
$ws.Range("X4").Value = 0.04621840177599682
$ws.Range("Y4").Value = 0.04601870691400484

$ws.Range("C5").Value = 0.7344110012054443
$ws.Range("E5").Value = 2202.126286495735
$ws.Range("F5").Value = 0.1497593480639168
$ws.Range("G5").Value = 0.1150776169237128
$ws.Range("H5").Value = 0.09420252304396427
$ws.Range("I5").Value = 0.08106652238228676
$ws.Range("J5").Value = 0.07032213466715975
$ws.Range("K5").Value = 0.06336161133801914
$ws.Range("L5").Value = 0.0587535044603494
$ws.Range("M5").Value = 0.05270396421690392
$ws.Range("N5").Value = 0.05127780498285568
$ws.Range("O5").Value = 0.05024495653787084
$ws.Range("P5").Value = 0.04865183969592821
$ws.Range("Q5").Value = 0.04715726074969269
$ws.Range("R5").Value = 0.04587877057471723
$ws.Range("S5").Value = 0.04544753643029661
$ws.Range("T5").Value = 0.04445845528026365
$ws.Range("U5").Value = 0.04425215192444547
$ws.Range("V5").Value = 0.04342089270836839
$ws.Range("W5").Value = 0.04333439651787391
$ws.Range("X5").Value = 0.0430182666002697
$ws.Range("Y5").Value = 0.04292643833325019

$ws.Range("C6").Value = 0.7343628406524658
$ws.Range("E6").Value = 2238.452593445997
$ws.Range("F6").Value = 0.1497593480639168
$ws.Range("G6").Value = 0.1273342249029688
$ws.Range("H6").Value = 0.09834562226649746
$ws.Range("I6").Value = 0.08278138951877775
$ws.Range("J6").Value = 0.07417545257611062
$ws.Range("K6").Value = 0.06611260427658407
$ws.Range("L6").Value = 0.06176315784337776
$ws.Range("M6").Value = 0.0573189936903312
$ws.Range("N6").Value = 0.05381713138375946
$ws.Range("O6").Value = 0.05022147756630944
$ws.Range("P6").Value = 0.04935866520371448
$ws.Range("Q6").Value = 0.04750184179756599
$ws.Range("R6").Value = 0.04723804836265769
$ws.Range("S6").Value = 0.0459666389517788
$ws.Range("T6").Value = 0.0451779414181981
$ws.Range("U6").Value = 0.04476021280251491
$ws.Range("V6").Value = 0.04421328817394798
$ws.Range("W6").Value = 0.04391281415093422
$ws.Range("X6").Value = 0.04376813959888919
$ws.Range("Y6").Value = 0.04363455347847947

$ws.Range("C7").Value = 0.7187647819519043
$ws.Range("E7").Value = 2322.55232054115
$ws.Range("F7").Value = 0.1497593480639168
$ws.Range("G7").Value = 0.121051373265077
$ws.Range("H7").Value = 0.1015462242439347
$ws.Range("I7").Value = 0.08785649886713838
$ws.Range("J7").Value = 0.07653748951497456
$ws.Range("K7").Value = 0.07108744693888748
$ws.Range("L7").Value = 0.06532285974153362
$ws.Range("M7").Value = 0.0623509121110928
$ws.Range("N7").Value = 0.05563673721776942
$ws.Range("O7").Value = 0.05378637043887775
$ws.Range("P7").Value = 0.0507922605342291
$ws.Range("Q7").Value = 0.0504351380058428
$ws.Range("R7").Value = 0.04888747650651677
$ws.Range("S7").Value = 0.04777419082765217
$ws.Range("T7").Value = 0.04691947740563192
$ws.Range("U7").Value = 0.04647352385885735
$ws.Range("V7").Value = 0.04598677325518256
$ws.Range("W7").Value = 0.0457135603675387
$ws.Range("X7").Value = 0.04542238432670259
$ws.Range("Y7").Value = 0.04527392437702046

$ws.Range("C8").Value = 0.70308518409729
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 2401.849469895791
$ws.Range("F8").Value = 0.1290701352274599
$ws.Range("G8").Value = 0.1270133153463208
$ws.Range("H8").Value = 0.1117232063016584
$ws.Range("I8").Value = 0.09665971620824125
$ws.Range("J8").Value = 0.08206766966181848
$ws.Range("K8").Value = 0.07189040343058704
$ws.Range("L8").Value = 0.06865720027486887
$ws.Range("M8").Value = 0.06333910381666959
$ws.Range("N8").Value = 0.05855306231114522
$ws.Range("O8").Value = 0.05664882094386005
$ws.Range("P8").Value = 0.05384359910907926
$ws.Range("Q8").Value = 0.05284422807656074
$ws.Range("R8").Value = 0.050990567363592
$ws.Range("S8").Value = 0.04995369644771817
$ws.Range("T8").Value = 0.04823994171586801
$ws.Range("U8").Value = 0.04797003321811097
$ws.Range("V8").Value = 0.04767266914647879
$ws.Range("W8").Value = 0.0472731101195637
$ws.Range("X8").Value = 0.04689723502356903
$ws.Range("Y8").Value = 0.04681967777574641

$ws.Range("C9").Value = 0.7187929153442383
$ws.Range("E9").Value = 2283.044177400936
$ws.Range("F9").Value = 0.1237898626793374
$ws.Range("G9").Value = 0.1237898626793374
$ws.Range("H9").Value = 0.1137535194410762
$ws.Range("I9").Value = 0.09044289117800083
$ws.Range("J9").Value = 0.07958117766702215
$ws.Range("K9").Value = 0.06630613126913937
$ws.Range("L9").Value = 0.0611658136712826
$ws.Range("M9").Value = 0.05890219355154275
$ws.Range("N9").Value = 0.05540752689693233
$ws.Range("O9").Value = 0.05230648169865215
$ws.Range("P9").Value = 0.05065512412086562
$ws.Range("Q9").Value = 0.04960146059780877
$ws.Range("R9").Value = 0.04778493808594001
$ws.Range("S9").Value = 0.04648785251838729
$ws.Range("T9").Value = 0.04646220983655502
$ws.Range("U9").Value = 0.04548149340623794
$ws.Range("V9").Value = 0.04517396984899325
$ws.Range("W9").Value = 0.04469789586475171
$ws.Range("X9").Value = 0.04457158786177876
$ws.Range("Y9").Value = 0.04450378513452116

$ws.Range("C10").Value = 0.8281145095825195
$ws.Range("E10").Value = 2310.025359150329
$ws.Range("F10").Value = 0.1497593480639168
$ws.Range("G10").Value = 0.1184801079333656
$ws.Range("H10").Value = 0.1007028535786967
$ws.Range("I10").Value = 0.08353091886654077
$ws.Range("J10").Value = 0.07562289750995319
$ws.Range("K10").Value = 0.06734268792914982
$ws.Range("L10").Value = 0.06188819093126791
$ws.Range("M10").Value = 0.06046636076935589
$ws.Range("N10").Value = 0.05618603323469928
$ws.Range("O10").Value = 0.0544932049798906
$ws.Range("P10").Value = 0.05198501318557715
$ws.Range("Q10").Value = 0.0503935885499279
$ws.Range("R10").Value = 0.04920241166147399
$ws.Range("S10").Value = 0.04748584229543912
$ws.Range("T10").Value = 0.04657422757736897
$ws.Range("U10").Value = 0.04608498473481625
$ws.Range("V10").Value = 0.04544081908422815
$ws.Range("W10").Value = 0.04544081908422815
$ws.Range("X10").Value = 0.04518872069346747
$ws.Range("Y10").Value = 0.04502973409649763

$ws.Range("C11").Value = 0.7499945163726807
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 2367.328310296081
$ws.Range("F11").Value = 0.1245793449342815
$ws.Range("G11").Value = 0.1202804491017233
$ws.Range("H11").Value = 0.1108205439279128
$ws.Range("I11").Value = 0.09383266358783419
$ws.Range("J11").Value = 0.080107417907318
$ws.Range("K11").Value = 0.07135732214029818
$ws.Range("L11").Value = 0.0636087220957266
$ws.Range("M11").Value = 0.06211569556048255
$ws.Range("N11").Value = 0.05815528875412735
$ws.Range("O11").Value = 0.05575510449337018
$ws.Range("P11").Value = 0.05313662278201187
$ws.Range("Q11").Value = 0.05062346508681892
$ws.Range("R11").Value = 0.04957443010477465
$ws.Range("S11").Value = 0.04908954757613705
$ws.Range("T11").Value = 0.04828367135735194
$ws.Range("U11").Value = 0.04750679189229829
$ws.Range("V11").Value = 0.04695450001941078
$ws.Range("W11").Value = 0.04648727717124724
$ws.Range("X11").Value = 0.04648727717124724
$ws.Range("Y11").Value = 0.04614675068803277
